$d   = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function XmlEscape($s) {
    if ($null -eq $s) { return "" }
    return $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

# ---------------------------------------------------------------------------
# Footer: the "NUMPAGES" field is currently stored as a legacy <w:fldSimple>
# run. Re-express it as an explicit begin/separate/end complex field (the
# same shape already used for the "PAGE" field right before it), which is
# what Word normally canonicalizes this kind of field to.
# ---------------------------------------------------------------------------
$ftr = $sec.Footers.Item(1)

# wdFieldPage = 33, wdFieldNumPages = 26 - locate by type rather than a
# fixed index so this keeps working regardless of field order.
$pageField = $null
$numField  = $null
for ($i = 1; $i -le $ftr.Range.Fields.Count; $i++) {
    $candidate = $ftr.Range.Fields.Item($i)
    if ($candidate.Type -eq 33) { $pageField = $candidate }
    if ($candidate.Type -eq 26) { $numField = $candidate }
}

$leadText = $ftr.Range.Duplicate
$leadText.Start = 0
$leadText.End = $pageField.Code.Start - 1

$midText = $ftr.Range.Duplicate
$midText.Start = $pageField.Result.End + 1
$midText.End = $numField.Code.Start - 1

$leadText = XmlEscape $leadText.Text
$midText = XmlEscape $midText.Text
$pageCode = XmlEscape $pageField.Code.Text
$pageResult = XmlEscape $pageField.Result.Text
$numCode = XmlEscape $numField.Code.Text
$numResult = XmlEscape $numField.Result.Text

$ftrBody = "<w:p><w:pPr><w:pStyle w:val=`"Footer`"/></w:pPr>" +
  "<w:r><w:t xml:space=`"preserve`">$leadText</w:t></w:r>" +
  "<w:r><w:fldChar w:fldCharType=`"begin`"/></w:r>" +
  "<w:r><w:instrText xml:space=`"preserve`">$pageCode</w:instrText></w:r>" +
  "<w:r><w:fldChar w:fldCharType=`"separate`"/></w:r>" +
  "<w:r><w:rPr><w:noProof/></w:rPr><w:t>$pageResult</w:t></w:r>" +
  "<w:r><w:fldChar w:fldCharType=`"end`"/></w:r>" +
  "<w:r><w:t xml:space=`"preserve`">$midText</w:t></w:r>" +
  "<w:r><w:fldChar w:fldCharType=`"begin`"/></w:r>" +
  "<w:r><w:instrText xml:space=`"preserve`">$numCode</w:instrText></w:r>" +
  "<w:r><w:fldChar w:fldCharType=`"separate`"/></w:r>" +
  "<w:r><w:rPr><w:noProof/></w:rPr><w:t>$numResult</w:t></w:r>" +
  "<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType=`"end`"/></w:r>" +
  "</w:p>"

$ftrXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  $ftrBody +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$ftrRange = $ftr.Range.Duplicate
[void]$ftrRange.InsertXML($ftrXml)

# ---------------------------------------------------------------------------
# First-page header: drop the stray empty run (<w:r><w:t/></w:r>) that adds
# no content - the paragraph should just keep its Header style.
# ---------------------------------------------------------------------------
$hdr2 = $sec.Headers.Item(2)

$hdrXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:pPr><w:pStyle w:val="Header"/></w:pPr></w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$hdr2Range = $hdr2.Range.Duplicate
[void]$hdr2Range.InsertXML($hdrXml)
